$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text cells (Coin name / Link) -- rows 9, 10, 51 swapped/updated
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'

# Force column D (Price) to remain plain text so values like "54.09" are not
# coerced to numbers -- matches the workbook author inlineStr/text convention.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '36.307.80'
$ws.Range('D3').Value = '1.971.86'
$ws.Range('D5').Value = '233.05'
$ws.Range('D8').Value = '54.09'
$ws.Range('D9').Value = '0.369'
$ws.Range('D10').Value = '58.34'
$ws.Range('D12').Value = '0.0983'
$ws.Range('D13').Value = '2.262.01'
$ws.Range('D14').Value = '13.85'
$ws.Range('D16').Value = '0.747'
$ws.Range('D17').Value = '5.01'
$ws.Range('D18').Value = '1.970.65'
$ws.Range('D19').Value = '36.283.16'
$ws.Range('D20').Value = '67.51'
$ws.Range('D21').Value = '0.0₃0804'
$ws.Range('D22').Value = '5.22'
$ws.Range('D23').Value = '220.79'
$ws.Range('D26').Value = '2.34'
$ws.Range('D27').Value = '160.21'
$ws.Range('D28').Value = '8.51'
$ws.Range('D29').Value = '18.66'
$ws.Range('D34').Value = '0.0601'
$ws.Range('D35').Value = '4.21'
$ws.Range('D39').Value = '3.22'
$ws.Range('D40').Value = '5.24'
$ws.Range('D41').Value = '3.03'
$ws.Range('D42').Value = '1.447.37'
$ws.Range('D43').Value = '0.0884'
$ws.Range('D45').Value = '87.97'
$ws.Range('D47').Value = '0.985'
$ws.Range('D48').Value = '14.70'
$ws.Range('D51').Value = '2.154.95'

# Restore the default cell style so no stray number-format style is left
# referenced on these cells (keeps styles.xml / cell format untouched).
$priceRange.Style = "Normal"

# Volume(1h) percentage cells -- plain text already (contains "%" and spaces).
$ws.Range('E2').Value = '  -3.33%  '
$ws.Range('E3').Value = '  -2.22%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  -11.63%  '
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -3.39%  '
$ws.Range('E9').Value = '  -4.43%  '
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('E11').Value = '  -3.52%  '
$ws.Range('E12').Value = '  -3.27%  '
$ws.Range('E13').Value = '  -2.05%  '
$ws.Range('E14').Value = '  -3.50%  '
$ws.Range('E15').Value = '  -4.92%  '
$ws.Range('E16').Value = '  -7.33%  '
$ws.Range('E17').Value = '  -4.65%  '
$ws.Range('E18').Value = '  -1.75%  '
$ws.Range('E19').Value = '  -3.15%  '
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('E21').Value = '  -4.65%  '
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('E23').Value = '  -3.35%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('E26').Value = '  -12.85%  '
$ws.Range('E27').Value = '  -2.58%  '
$ws.Range('E28').Value = '  -4.98%  '
$ws.Range('E29').Value = '  -5.13%  '
$ws.Range('E30').Value = '  -3.75%  '
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('E33').Value = '  -6.85%  '
$ws.Range('E34').Value = '  -7.68%  '
$ws.Range('E35').Value = '  -6.98%  '
$ws.Range('E36').Value = '  -4.39%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('E39').Value = '  -3.75%  '
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('E41').Value = '  -1.50%  '
$ws.Range('E42').Value = '  +3.50%  '
$ws.Range('E43').Value = '  -5.60%  '
$ws.Range('E44').Value = '  -6.22%  '
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('E46').Value = '  -11.80%  '
$ws.Range('E47').Value = '  -4.19%  '
$ws.Range('E48').Value = '  -6.17%  '
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('E50').Value = '  -4.92%  '
$ws.Range('E51').Value = '  -2.14%  '
